# Updates loading_percent.xlsx (Sheet1) values for the 380 kV case
# (Case_3_150/res_line/loading_percent.xlsx) across rows 2-25,
# columns B, C, D, F, G, L, O.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @{ "B"=17.34599836510662; "C"=11.80355612291052; "D"=4.763576408364015; "F"=24.07637840253037; "G"=3.626892908436322; "L"=10.68693129779342; "O"=21.49639139889986 }
  3 = @{ "B"=16.68438401551996; "C"=11.62877742975136; "D"=4.730888068236987; "F"=24.11704897371636; "G"=3.629230137562073; "L"=10.6605434325542; "O"=21.59038533749178 }
  4 = @{ "B"=16.26584124609315; "C"=11.52039494055212; "D"=4.710693224892956; "F"=24.1512997303878; "G"=3.630740728564255; "L"=10.64654539395664; "O"=21.65520489665389 }
  5 = @{ "B"=16.0924294813344; "C"=11.47600134922573; "D"=4.702436592914898; "F"=24.1675792204104; "G"=3.631375359955832; "L"=10.64139923842723; "O"=21.68339677888074 }
  6 = @{ "B"=16.06346995360765; "C"=11.46861735782336; "D"=4.701064070915101; "F"=24.17042230371456; "G"=3.631481892624139; "L"=10.64057853803337; "O"=21.68818508016698 }
  7 = @{ "B"=16.26351376199495; "C"=11.51979709599043; "D"=4.71058197723869; "F"=24.15150989551022; "G"=3.630749210195991; "L"=10.64647372638968; "O"=21.6555779187436 }
  8 = @{ "B"=17.12056599841001; "C"=11.74353873623275; "D"=4.752333697104516; "F"=24.08846951891266; "G"=3.62768314792349; "L"=10.67737776486871; "O"=21.5273197219727 }
  9 = @{ "B"=18.69422880829824; "C"=12.1720427316766; "D"=4.833043863193208; "F"=24.03889744482203; "G"=3.622266996460513; "L"=10.75527663621925; "O"=21.33261671540938 }
  10 = @{ "B"=19.77437357126453; "C"=12.47830880223908; "D"=4.89139183997344; "F"=24.04808907893902; "G"=3.618647326835275; "L"=10.82276221509853; "O"=21.22475710247085 }
  11 = @{ "B"=20.24742276612968; "C"=12.61531707549802; "D"=4.917679031752522; "F"=24.0622372531042; "G"=3.617077872730621; "L"=10.85561723277127; "O"=21.18344135714538 }
  12 = @{ "B"=20.42379472742818; "C"=12.66682918733896; "D"=4.927592296833291; "F"=24.06903004278417; "G"=3.616494590948344; "L"=10.86836164481043; "O"=21.16891915910962 }
  13 = @{ "B"=20.3859345394798; "C"=12.65575216088796; "D"=4.925459195887486; "F"=24.06750325371124; "G"=3.616619721186104; "L"=10.86560355140582; "O"=21.17199668672945 }
  14 = @{ "B"=20.26198900856468; "C"=12.61956262955325; "D"=4.918495451082054; "F"=24.06276733261129; "G"=3.617029664931882; "L"=10.85665969389423; "O"=21.18222404871437 }
  15 = @{ "B"=20.18570566672028; "C"=12.59734624205353; "D"=4.914224471883544; "F"=24.06005337176236; "G"=3.617282202868879; "L"=10.85122055712126; "O"=21.18863512507739 }
  16 = @{ "B"=19.74307882092131; "C"=12.46930521018738; "D"=4.88966839236384; "F"=24.04736523397434; "G"=3.61875144188363; "L"=10.8206578460571; "O"=21.22761389174031 }
  17 = @{ "B"=19.46675412766277; "C"=12.39013656587294; "D"=4.87453544573854; "F"=24.04213624626151; "G"=3.61967249152227; "L"=10.80245583056448; "O"=21.25351744408405 }
  18 = @{ "B"=19.30610211009949; "C"=12.34438520344726; "D"=4.865807576019749; "F"=24.04006680562476; "G"=3.620209520237216; "L"=10.79218984768787; "O"=21.26914546772143 }
  19 = @{ "B"=19.25141733233956; "C"=12.32885870183657; "D"=4.862848515673679; "F"=24.03952716201094; "G"=3.620392598530905; "L"=10.78874909191471; "O"=21.27456177895334 }
  20 = @{ "B"=19.49634808690645; "C"=12.39858679747897; "D"=4.876148868643495; "F"=24.04259576695457; "G"=3.619573692702173; "L"=10.80437247124481; "O"=21.25068445481391 }
  21 = @{ "B"=20.29847069441304; "C"=12.63020270084493; "D"=4.920542022413593; "F"=24.06411943118408; "G"=3.61690895544986; "L"=10.85927855843891; "O"=21.17918947465659 }
  22 = @{ "B"=20.80656072651161; "C"=12.77940467632309; "D"=4.94931382161916; "F"=24.08655123992911; "G"=3.615231697757849; "L"=10.89692499782763; "O"=21.13901382617841 }
  23 = @{ "B"=20.53689845052159; "C"=12.6999834128005; "D"=4.933981316650288; "F"=24.07381344298359; "G"=3.616121017133286; "L"=10.87667358477604; "O"=21.15985419915949 }
  24 = @{ "B"=19.4829742212286; "C"=12.39476718081056; "D"=4.875419526193646; "F"=24.04238509978239; "G"=3.619618336266205; "L"=10.80350533885797; "O"=21.2519629581249 }
  25 = @{ "B"=18.281172015905; "C"=12.05746162088409; "D"=4.811359109535113; "F"=24.04432491042455; "G"=3.623668776202825; "L"=10.73237926490968; "O"=21.37914832880835 }
}

foreach ($row in $data.Keys) {
  foreach ($col in $data[$row].Keys) {
    $addr = "$col$row"
    $ws.Range($addr).Value = $data[$row][$col]
  }
}
